$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used after 'Jofra Archer' (matches source data)
$nbsp = [char]0x00A0
$batsman = 'Jofra Archer' + $nbsp

# New match rows appended to the table (rows 12-21)
$rows = @(
    @(12, ' Dubai (DSC)', ' October 14 2020', 'Capitals won by 13 runs', 'Rajasthan Royals', 'Delhi Capitals', '1', '4', '0', '0', '25.00'),
    @(13, ' Dubai (DSC)', ' October 22 2020', 'Sunrisers won by 8 wickets (with 11 balls remaining)', 'Rajasthan Royals', 'Sunrisers Hyderabad', '16', '7', '1', '1', '228.57'),
    @(14, ' Sharjah', ' September 22 2020', 'Royals won by 16 runs', 'Rajasthan Royals', 'Chennai Super Kings', '27', '8', '0', '4', '337.50'),
    @(15, ' Dubai (DSC)', ' September 30 2020', 'KKR won by 37 runs', 'Rajasthan Royals', 'Kolkata Knight Riders', '6', '4', '0', '1', '150.00'),
    @(16, ' Dubai (DSC)', ' November 01 2020', 'KKR won by 60 runs', 'Rajasthan Royals', 'Kolkata Knight Riders', '6', '9', '0', '0', '66.66'),
    @(17, ' Sharjah', ' September 27 2020', 'Royals won by 4 wickets (with 3 balls remaining)', 'Rajasthan Royals', 'Kings XI Punjab', '13', '3', '0', '2', '433.33'),
    @(18, ' Dubai (DSC)', ' October 17 2020', 'RCB won by 7 wickets (with 2 balls remaining)', 'Rajasthan Royals', 'Royal Challengers Bangalore', '2', '3', '0', '0', '66.66'),
    @(19, ' Abu Dhabi', ' October 06 2020', 'Mumbai won by 57 runs', 'Rajasthan Royals', 'Mumbai Indians', '24', '11', '3', '1', '218.18'),
    @(20, ' Abu Dhabi', ' October 03 2020', 'RCB won by 8 wickets (with 5 balls remaining)', 'Rajasthan Royals', 'Royal Challengers Bangalore', '16', '10', '1', '1', '160.00'),
    @(21, ' Sharjah', ' October 09 2020', 'Capitals won by 46 runs', 'Rajasthan Royals', 'Delhi Capitals', '2', '4', '0', '0', '50.00'),
)

foreach ($row in $rows) {
    $r = $row[0]
    $venue = $row[1]
    $date = $row[2]
    $result = $row[3]
    $ownTeam = $row[4]
    $oppTeam = $row[5]
    $totalRuns = $row[6]
    $totalBalls = $row[7]
    $total4s = $row[8]
    $total6s = $row[9]
    $sr = $row[10]

    $rowRange = $ws.Range("A" + $r + ":K" + $r)
    $rowRange.NumberFormat = "@"

    $ws.Range("A" + $r).Value = $venue
    $ws.Range("B" + $r).Value = $date
    $ws.Range("C" + $r).Value = $result
    $ws.Range("D" + $r).Value = $ownTeam
    $ws.Range("E" + $r).Value = $oppTeam
    $ws.Range("F" + $r).Value = $batsman
    $ws.Range("G" + $r).Value = $totalRuns
    $ws.Range("H" + $r).Value = $totalBalls
    $ws.Range("I" + $r).Value = $total4s
    $ws.Range("J" + $r).Value = $total6s
    $ws.Range("K" + $r).Value = $sr
}

# The whole table holds text that looks numeric (as in the original rows),
# so keep Excel's "Number Stored as Text" warning suppressed across the
# full, now-larger range (originally only A1:K11 was marked this way).
$fullRange = $ws.Range("A1:K21")
try {
    $fullRange.Errors.Item(9).Ignore = $true
} catch {
}

